$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '56.684.20'
$ws.Range('E2').Value = '  +3.82%  '

$ws.Range('D3').Value = '3.248.17'
$ws.Range('E3').Value = '  +2.00%  '

$ws.Range('E4').Value = '  +0.02%  '

$ws.Range('D5').Value = "'395.59"
$ws.Range('E5').Value = '  -1.69%  '

$ws.Range('D6').Value = "'108.76"
$ws.Range('E6').Value = '  +0.06%  '

$ws.Range('D7').Value = "'0.583"
$ws.Range('E7').Value = '  +5.84%  '

$ws.Range('D8').Value = '3.244.52'
$ws.Range('E8').Value = '  +1.98%  '

$ws.Range('D9').Value = "'0.999"

$ws.Range('E10').Value = '  +1.01%  '

$ws.Range('D11').Value = "'39.21"
$ws.Range('E11').Value = '  +0.27%  '

$ws.Range('D12').Value = "'0.0976"
$ws.Range('E12').Value = '  +10.13%  '

$ws.Range('E13').Value = '  +2.05%  '

$ws.Range('D14').Value = '3.761.02'
$ws.Range('E14').Value = '  +2.20%  '

$ws.Range('E15').Value = '  +3.40%  '

$ws.Range('E16').Value = '  +0.37%  '

$ws.Range('D17').Value = '3.260.09'
$ws.Range('E17').Value = '  +2.43%  '

$ws.Range('E18').Value = '  -3.05%  '

$ws.Range('D19').Value = "'10.74"
$ws.Range('E19').Value = '  +1.75%  '

$ws.Range('D20').Value = '56.689.64'
$ws.Range('E20').Value = '  +3.93%  '

$ws.Range('E21').Value = '  +0.90%  '

$ws.Range('E22').Value = '  +8.25%  '

$ws.Range('D23').Value = "'12.92"
$ws.Range('E23').Value = '  -0.14%  '

$ws.Range('D24').Value = "'297.57"
$ws.Range('E24').Value = '  +8.11%  '

$ws.Range('E25').Value = '  +2.02%  '

$ws.Range('E26').Value = '  -3.29%  '

$ws.Range('D27').Value = "'28.12"
$ws.Range('E27').Value = '  +0.84%  '

$ws.Range('E28').Value = '  +0.59%  '

$ws.Range('D29').Value = "'7.66"
$ws.Range('E29').Value = '  -4.62%  '

$ws.Range('D30').Value = "'7.33"
$ws.Range('E30').Value = '  -2.76%  '

$ws.Range('E31').Value = '  -0.85%  '

$ws.Range('E32').Value = '  +0.05%  '

$ws.Range('D33').Value = "'11.23"
$ws.Range('E33').Value = '  +1.66%  '

$ws.Range('E34').Value = '  -3.50%  '

$ws.Range('D35').Value = "'39.65"
$ws.Range('E35').Value = '  +7.39%  '

$ws.Range('D36').Value = "'0.0487"
$ws.Range('E36').Value = '  -3.30%  '

$ws.Range('E37').Value = '  +2.03%  '

$ws.Range('D38').Value = "'51.35"
$ws.Range('E38').Value = '  +0.55%  '

$ws.Range('E39').Value = '  -0.08%  '

$ws.Range('D40').Value = "'3.48"
$ws.Range('E40').Value = '  -4.64%  '

$ws.Range('E41').Value = '  +1.18%  '

$ws.Range('D42').Value = "'137.98"
$ws.Range('E42').Value = '  +5.29%  '

$ws.Range('E43').Value = '  +3.70%  '

$ws.Range('E44').Value = '  -4.48%  '

$ws.Range('E45').Value = '  -2.82%  '

$ws.Range('E46').Value = '  -1.43%  '

$ws.Range('E47').Value = '  -3.55%  '

$ws.Range('D48').Value = "'22.27"
$ws.Range('E48').Value = '  -0.25%  '

$ws.Range('E49').Value = '  +3.50%  '

$ws.Range('D50').Value = '2.156.55'
$ws.Range('E50').Value = '  +3.11%  '

$ws.Range('E51').Value = '  -5.38%  '

